$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# GPA column (E) holds text-formatted numbers like "3.88" in the source data,
# not real numbers - force text number format so values aren't coerced.
$ws.Range("E2:E7").NumberFormat = "@"

# Row 2: Math -> Calculus, GPA 4.00 -> 3.88
$ws.Range("B2").Value = "Calculus"
$ws.Range("E2").Value = "3.88"

# Row 3: B3 "2" -> "English", GPA 4.00 -> 3.88 (C3 credit stays 2)
$ws.Range("B3").Value = "English"
$ws.Range("E3").Value = "3.88"

# Row 4: Semester 2 -> Semester 1, Englishj -> PROBLEM SOLVING AND PROGRAMMING,
# Grade B -> A-, GPA 3.00 -> 3.88
$ws.Range("A4").Value = "Semester 1"
$ws.Range("B4").Value = "PROBLEM SOLVING AND PROGRAMMING"
$ws.Range("D4").Value = "A-"
$ws.Range("E4").Value = "3.88"

# New row 5: Semester 2, FUNDAMENTALS OF COMPUTER NETWORKS (leading tab), 4, B, 3.10
$ws.Range("A5").Value = "Semester 2"
$ws.Range("B5").Value = "`tFUNDAMENTALS OF COMPUTER NETWORKS"
$ws.Range("C5").Value = 4
$ws.Range("D5").Value = "B"
$ws.Range("E5").Value = "3.10"

# New row 6: Semester 2, PROBABILITY AND STATISTICS, 3, B+, 3.10
$ws.Range("A6").Value = "Semester 2"
$ws.Range("B6").Value = "PROBABILITY AND STATISTICS"
$ws.Range("C6").Value = 3
$ws.Range("D6").Value = "B+"
$ws.Range("E6").Value = "3.10"

# New row 7: Semester 2, COMPUTER ARCHITECTURE, 3, B, 3.10
$ws.Range("A7").Value = "Semester 2"
$ws.Range("B7").Value = "COMPUTER ARCHITECTURE"
$ws.Range("C7").Value = 3
$ws.Range("D7").Value = "B"
$ws.Range("E7").Value = "3.10"
